$wb = $excel.ActiveWorkbook

# --- "About" sheet updates -------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Switch from U.S. units to EU units
$wsAbout.Range("A10").Value = "For the EU.:"
$wsAbout.Range("A11").Value = "The large primary energy output unit (used in totals graphs) is: TWh"
$wsAbout.Range("A12").Value = "The small primary energy output unit (used in energy intensity per unit GDP graphs) is: MWh"

# New "Relevant Conversion Factors" block
$wsAbout.Range("A15").Value = "Relevant Conversion Factors"
$wsAbout.Range("A15").Font.Bold = $true
$wsAbout.Range("A15:B15").Interior.Color = 12632256

$wsAbout.Range("A16").Value = 3412000
$wsAbout.Range("B16").Value = "BTU/MWh (pure unit conversion, not a heat rate)"

# Column widths (~13.29 / ~15.29 characters)
$wsAbout.Columns("A").ColumnWidth = 12.45
$wsAbout.Columns("B").ColumnWidth = 14.45

# --- "BpTPEU-large" sheet: now derives from About!A16 ----------------------
$wsLarge = $wb.Worksheets.Item("BpTPEU-large")
$wsLarge.Range("B2").Formula = "=About!A16*10^6"
$wsLarge.Range("B2").Select()

# --- "BpTPEU-small" sheet: now derives from About!A16, default style -------
$wsSmall = $wb.Worksheets.Item("BpTPEU-small")
$wsSmall.Range("B2").Style = "Normal"
$wsSmall.Range("B2").Formula = "=About!A16"

# Make BpTPEU-small the active/selected sheet, mirroring the author's last
# saved view.
$wsSmall.Activate()
$wsSmall.Range("F23").Select()
